$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.004.29'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.867.08'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9983'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5094'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3875'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08324'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.113'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.204'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').Value = '1.842.81'
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.222'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9997'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001098'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06658'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9972'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.962'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = '28.006.84'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.240'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.457'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.38%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '125.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1049'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.030'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.777'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.580'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('B34').Value = 'FraxShare'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.500'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02431'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.90%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06518'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2199'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.189'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6438'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.222'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.939'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.19'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'Decentraland'
$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6069'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.79%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.271'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.640'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.997'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.228'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '119.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06872'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.31'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.40%  '
